$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.577.18"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.908.26"
$ws.Range("E3").Value = "  +4.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.83"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.86"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("D7").Value = "3.906.25"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -3.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.97"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "4.561.38"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "3.930.85"
$ws.Range("E16").Value = "  +5.07%  "
$ws.Range("D17").Value = "68.760.50"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.42"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.18"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "484.86"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.718"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("E24").Value = "  +13.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.35"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.05"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "4.057.17"
$ws.Range("E31").Value = "  +4.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.88"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.91"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "3.850.85"
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.87"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.319"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "439.76"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.46"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "2.835.48"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.92"
$ws.Range("E50").Value = "  +10.28%  "
$ws.Range("E51").Value = "  +0.56%  "